$d = $word.ActiveDocument

# Locate the unique merge-field result text "«=date.first»". It lives in its
# own paragraph inside a table cell, immediately followed by the paragraph
# that closes the "contact.dates_by_medium_type" merge loop
# ("«contact.dates_by_medium_type:endEach»"). We need to insert a new, empty
# paragraph between those two paragraphs (for better multi-line formatting in
# the table).
$hitRange = $d.Content
$found = $hitRange.Find.Execute("«=date.first»", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate anchor text for the edit"
}

# Resolve the actual paragraph that contains the found text via the
# document's Paragraphs collection (more reliable here than Range.Paragraphs
# on a short/zero-length range).
$count = $d.Paragraphs.Count
$target = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($hitRange.Start -ge $p.Range.Start -and $hitRange.Start -lt $p.Range.End) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not resolve the paragraph containing the anchor text"
}

$paraEnd = $target.Range.End

# Build a minimal, self-contained WordprocessingML package describing the new
# (empty) paragraph, matching the paragraph-mark formatting used throughout
# this table (Times New Roman complex-script font, size 12pt/24 half-points).
# Inserting via InsertXML avoids the engine's paragraph-split formatting
# heuristics (which can pick up unrelated direct formatting, e.g. from the
# document title) producing an exact, clean result.
$newParagraphXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:rFonts w:cs="Times New Roman"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertPoint = $d.Range($paraEnd, $paraEnd)
$insertPoint.InsertXML($newParagraphXml) | Out-Null
